# Update the version number shown in the title of the cheat sheet
# ("nctoolkit v0.9.3" -> "nctoolkit v1.1.6") on slide 1, without
# disturbing any other runs/formatting in the title text box.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)
$sh = $s.Shapes.Item(1)
$tr = $sh.TextFrame.TextRange

$oldVersion = " v0.9.3"
$newVersion = " v1.1.6"

$fullText = $tr.Text
$idx = $fullText.IndexOf($oldVersion)
if ($idx -ge 0) {
    $run = $tr.Characters($idx + 1, $oldVersion.Length)
    $run.Text = $newVersion
}
